$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 43317
$ws.Range("B2").Value = "Lucas Gabriel Nascimento"
$ws.Range("C2").Value = "Engenharia"
$ws.Range("D2").Value = "Outros"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 45081
$ws.Range("G2").Value = 6495.35

# Row 3
$ws.Range("A3").Value = 8840
$ws.Range("B3").Value = "Lívia das Neves"
$ws.Range("C3").Value = "Engenharia"
$ws.Range("D3").Value = "Consulta médica"
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 45104
$ws.Range("G3").Value = 5680.91

# Row 4
$ws.Range("A4").Value = 7129
$ws.Range("B4").Value = "Laura Almeida"
$ws.Range("C4").Value = "Engenharia"
$ws.Range("D4").Value = "Viagem de negócios"
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 45081
$ws.Range("G4").Value = 8599.049999999999

# Row 5
$ws.Range("A5").Value = 25414
$ws.Range("B5").Value = "Erick Silveira"
$ws.Range("C5").Value = "Financeiro"
$ws.Range("D5").Value = "Problemas pessoais"
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 45096
$ws.Range("G5").Value = 6678.92

# Row 6
$ws.Range("A6").Value = 1532
$ws.Range("B6").Value = "Eloah Gomes"
$ws.Range("C6").Value = "P&D"
$ws.Range("D6").Value = "Problemas pessoais"
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 45104
$ws.Range("G6").Value = 3789.14

# Row 7
$ws.Range("A7").Value = 52808
$ws.Range("B7").Value = "Maria da Cruz"
$ws.Range("C7").Value = "Recursos Humanos"
$ws.Range("D7").Value = "Consulta médica"
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 45078
$ws.Range("G7").Value = 8160.54

# Row 8
$ws.Range("A8").Value = 74576
$ws.Range("B8").Value = "Srta. Ana Laura Cardoso"
$ws.Range("C8").Value = "Recursos Humanos"
$ws.Range("D8").Value = "Consulta médica"
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 45092
$ws.Range("G8").Value = 8330.48

# Row 9
$ws.Range("A9").Value = 62081
$ws.Range("B9").Value = "Joana Castro"
$ws.Range("C9").Value = "Jurídico"
$ws.Range("D9").Value = "Problemas pessoais"
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 45078
$ws.Range("G9").Value = 7037.28

# Row 10
$ws.Range("A10").Value = 85638
$ws.Range("B10").Value = "Bianca Duarte"
$ws.Range("C10").Value = "Atendimento ao Cliente"
$ws.Range("D10").Value = "Problemas pessoais"
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = 45086
$ws.Range("G10").Value = 8120.77

# Row 11
$ws.Range("A11").Value = 8243
$ws.Range("B11").Value = "Gustavo Henrique das Neves"
$ws.Range("C11").Value = "Jurídico"
$ws.Range("D11").Value = "Outros"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 45084
$ws.Range("G11").Value = 11431.23
